$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Structural changes -------------------------------------------------
# Insert a new header row at the very top. This shifts the existing rows
# (old A1:A7) down to A2:A8.
$ws.Rows("1:1").Insert()

# Insert two new rows right after the (now shifted) "Leg Motion Study
# Animation" row (row 8), to host the two new Design Tasks entries.
$ws.Rows("9:10").Insert()

# Row 11 is intentionally left blank (separator before the new "Report
# Tasks" section), then 5 new rows are inserted starting at row 12 for
# the "Report Tasks" section (1 header + 4 data rows).
$ws.Rows("12:16").Insert()

# --- New header row (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Task"
$ws.Range("B1").Value = "Taskee"
$ws.Range("C1").Value = "Definition of Finished"

# Copy the header formatting (italic, no wrap / italic+wrap for column C)
# from the existing "People"/"Definition of Finished" header used on the
# other weekly sheets so the new header matches the established style.
$src = $wb.Worksheets.Item(3).Range("A2:C2")
$src.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Fill in People / Definition of Finished for existing rows 3-8 ------
# (Row 2 "Design Tasks" is a section header, spanning only column A.)
$ws.Range("B3").Value = "Tyler"
$ws.Range("C3").Value = "Electrical components sized and relative positions mapped out"

$ws.Range("B4").Value = "Tyler"
$ws.Range("C4").Value = "Initial components known and placed for signal conditioning"

$ws.Range("B5").Value = "Tyler"
$ws.Range("C5").Value = "Determine style of Op Amp and number of amplification voltages"

$ws.Range("B6").Value = "Ron"
$ws.Range("C6").Value = "2 Shoulder designs completed"

$ws.Range("B7").Value = "Ron"
$ws.Range("C7").Value = "2 Leg designs completed"

$ws.Range("B8").Value = "Logan"
$ws.Range("C8").Value = "Animation has realistic gait and determines min/max values"

# --- New Design Tasks rows (9-10) ---------------------------------------
$ws.Range("A9").Value = "Update SolidWorks Model of Robot"
$ws.Range("B9").Value = "Justin"
$ws.Range("C9").Value = "Robot updated with T-slotted aluminum chassis"

$ws.Range("A10").Value = "Lightweight materials research"
$ws.Range("B10").Value = "Justin"
$ws.Range("C10").Value = "List of potential leg materials/weights"

# Give the new data rows (9-10) the same formatting as the rows above them.
$fmtSrc = $ws.Range("A8:C8")
$fmtSrc.Copy()
$ws.Range("A9:C10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 11 stays empty (blank separator row).

# --- New "Report Tasks" section (rows 12-16) ----------------------------
$ws.Range("A12").Value = "Report Tasks"
$hdrSrc = $ws.Range("A2")
$hdrSrc.Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A13").Value = "Background Section"
$ws.Range("B13").Value = "Justin"
$ws.Range("C13").Value = "Updates to background info, Kevin Lee moved"

$ws.Range("A14").Value = "Pedagogical Research"
$ws.Range("B14").Value = "Logan"
$ws.Range("C14").Value = "Updates to pedagogical research"

$ws.Range("A15").Value = "Motherboard/Electrical Diagram"
$ws.Range("B15").Value = "Tyler"
$ws.Range("C15").Value = "Clean up electrical section and add full wiring diagram"

$ws.Range("A16").Value = "Update pneumait diagram"
$ws.Range("B16").Value = "Logan"
$ws.Range("C16").Value = "Pneumatic diagram does not have an accumulation tank"

$dataFmtSrc = $ws.Range("A8:C8")
$dataFmtSrc.Copy()
$ws.Range("A13:C16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Update selection to mirror the authored workbook -------------------
$ws.Range("A12").Select()
